$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new hours entry row (row 12)
$ws.Range("A11").Copy()
$ws.Range("A12").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A12").Value = Get-Date -Year 2018 -Month 2 -Day 20 -Hour 0 -Minute 0 -Second 0 -Millisecond 0
$ws.Range("B12").Value = 2.5
$ws.Range("C12").Value = "weekly meeting; updating power function; started DataCamp on ggplots"

$ws.Range("A13").Select()
